$wb = $excel.ActiveWorkbook

# Turn off alerts so sheet deletion doesn't prompt for confirmation
$excel.DisplayAlerts = $false

# Delete the two empty placeholder sheets ("Sheet2" and "Sheet3"),
# keeping the data sheet ("Sheet4") and the first empty sheet ("Sheet1").
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Rename the data sheet from "Sheet4" to "Sheet2"
$wb.Worksheets.Item("Sheet4").Name = "Sheet2"

$excel.DisplayAlerts = $true

# Reset the selection on the data sheet back to A1 (was H2)
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()
$ws.Range("A1").Select()
